# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> linked from the notes master, currently "Office Theme"
#   ppt/theme/theme2.xml -> linked from the slide master (and the presentation
#                            itself), currently "Integral"
# The authored change swaps the two themes' content: theme1.xml ends up
# holding the "Integral" color scheme and theme2.xml ends up holding the
# "Office Theme" color scheme (the relationships/filenames themselves are
# untouched). Concretely this means: the slide master's theme colors become
# the stock "Office Theme" palette, and the notes master's theme colors
# become the "Integral" palette. Font scheme / format scheme are identical
# between the two themes already, so only the 12 color-scheme slots (and
# nothing else) need to change.

$p = $ppt.ActivePresentation

# RGB() style (r + g*256 + b*65536) encodings of each theme's 12 scheme
# colors, in the standard dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink order.
$officeThemeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)
$integralThemeColors = @(0, 16777215, 5332805, 13754083, 3722137, 3646819, 2412774, 38860, 13611854, 10915127, 2465643, 158642)

function Set-ThemeColorScheme($theme, $colors) {
    $scheme = $theme.ThemeColorScheme
    for ($i = 1; $i -le $colors.Length; $i++) {
        $scheme.Item($i).RGB = $colors[$i - 1]
    }
}

# Slide master's theme (backing file ppt/theme/theme2.xml) -> "Office Theme" colors
Set-ThemeColorScheme $p.SlideMaster.Theme $officeThemeColors

# Notes master's theme (backing file ppt/theme/theme1.xml) -> "Integral" colors
Set-ThemeColorScheme $p.NotesMaster.Theme $integralThemeColors
